$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($ws, $row1, $row2) {
    # Swap the contents of columns B through AD (columns 2-30) between two rows.
    # Column A (row index) and any unchanged columns are left untouched.
    # NOTE: Range.Value is a parameterized COM property and does not marshal
    # cleanly through this host, so use Value2 instead.
    $range1 = $ws.Range("B$row1" + ":AD$row1")
    $range2 = $ws.Range("B$row2" + ":AD$row2")

    $vals1 = $range1.Value2
    $vals2 = $range2.Value2

    $range1.Value2 = $vals2
    $range2.Value2 = $vals1
}

Swap-Rows $ws 58 59
Swap-Rows $ws 136 137
Swap-Rows $ws 141 142
